$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.098.58"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.374.00"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'302.83"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'97.01"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").Value = "'0.505"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "'18.32"
$ws.Range("E13").Value = "  -4.33%  "
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "2.747.51"
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").Value = "2.436.38"
$ws.Range("E16").Value = "  +6.04%  "
$ws.Range("E17").Value = "  +3.50%  "
$ws.Range("D18").Value = "43.076.79"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'12.18"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "'235.36"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "'24.79"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "'31.31"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'5.08"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'0.0751"
$ws.Range("E33").Value = "  +8.35%  "
$ws.Range("D34").Value = "'17.39"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").Value = "'0.104"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("E36").Value = "  +5.58%  "
$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "'4.28"
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("D39").Value = "'2.80"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").Value = "'22.31"
$ws.Range("E40").Value = "  +10.75%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'104.63"
$ws.Range("E42").Value = "  -36.61%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.957.71"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'9.16"
$ws.Range("E47").Value = "  -10.75%  "
$ws.Range("D48").Value = "2.602.39"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").Value = "'52.58"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  +1.76%  "
